# Add a new test-case row to the "Source" sheet:
# row 6 / column B gets a date value (43822 == 2019-12-23), formatted the
# same way as the existing B3/B4 date cells, and the active selection
# moves to D8 (matching the author's saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy B4's formatting (numFmtId 14 date format) onto B6, then set its value.
$ws.Range("B4").Copy($ws.Range("B6"))
$ws.Range("B6").Value = 43822

# Match the saved cursor/selection from the authored workbook.
$ws.Range("D8").Select()
